{"js": "const replacements = [\n  [\"821\\u00d74=\", \"580\\u00d77=\"],\n  [\"741\\u00d79=\", \"691\\u00d78=\"],\n  [\"177\\u00d78=\", \"767\\u00d74=\"],\n  [\"413\\u00d77=\", \"938\\u00d73=\"],\n  [\"645\\u00d79=\", \"749\\u00d78=\"],\n  [\"864\\u00d78=\", \"689\\u00d74=\"],\n  [\"796\\u00d77=\", \"264\\u00d75=\"],\n  [\"602\\u00d76=\", \"676\\u00d77=\"],\n  [\"547\\u00d73=\", \"555\\u00d77=\"],\n  [\"265\\u00d78=\", \"584\\u00d77=\"],\n  [\"475\\u00d79=\", \"929\\u00d74=\"],\n  [\"751\\u00d79=\", \"854\\u00d79=\"],\n  [\"454\\u00d79=\", \"893\\u00d74=\"],\n  [\"735\\u00d75=\", \"913\\u00d74=\"],\n  [\"785\\u00d78=\", \"490\\u00d78=\"],\n  [\"966\\u00d78=\", \"147\\u00d79=\"],\n  [\"915\\u00d77=\", \"430\\u00d72=\"],\n  [\"169\\u00d78=\", \"689\\u00d74=\"],\n  [\"344\\u00d77=\", \"328\\u00d77=\"],\n  [\"200\\u00d74=\", \"508\\u00d73=\"],\n  [\"925\\u00d73=\", \"941\\u00d78=\"],\n  [\"692\\u00d75=\", \"286\\u00d75=\"],\n  [\"252\\u00d74=\", \"588\\u00d72=\"],\n  [\"741\\u00d73=\", \"695\\u00d77=\"],\n  [\"262\\u00d77=\", \"619\\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"821\u00d74=\"; new = \"580\u00d77=\"},\n    @{old = \"741\u00d79=\"; new = \"691\u00d78=\"},\n    @{old = \"177\u00d78=\"; new = \"767\u00d74=\"},\n    @{old = \"413\u00d77=\"; new = \"938\u00d73=\"},\n    @{old = \"645\u00d79=\"; new = \"749\u00d78=\"},\n    @{old = \"864\u00d78=\"; new = \"689\u00d74=\"},\n    @{old = \"796\u00d77=\"; new = \"264\u00d75=\"},\n    @{old = \"602\u00d76=\"; new = \"676\u00d77=\"},\n    @{old = \"547\u00d73=\"; new = \"555\u00d77=\"},\n    @{old = \"265\u00d78=\"; new = \"584\u00d77=\"},\n    @{old = \"475\u00d79=\"; new = \"929\u00d74=\"},\n    @{old = \"751\u00d79=\"; new = \"854\u00d79=\"},\n    @{old = \"454\u00d79=\"; new = \"893\u00d74=\"},\n    @{old = \"735\u00d75=\"; new = \"913\u00d74=\"},\n    @{old = \"785\u00d78=\"; new = \"490\u00d78=\"},\n    @{old = \"966\u00d78=\"; new = \"147\u00d79=\"},\n    @{old = \"915\u00d77=\"; new = \"430\u00d72=\"},\n    @{old = \"169\u00d78=\"; new = \"689\u00d74=\"},\n    @{old = \"344\u00d77=\"; new = \"328\u00d77=\"},\n    @{old = \"200\u00d74=\"; new = \"508\u00d73=\"},\n    @{old = \"925\u00d73=\"; new = \"941\u00d78=\"},\n    @{old = \"692\u00d75=\"; new = \"286\u00d75=\"},\n    @{old = \"252\u00d74=\"; new = \"588\u00d72=\"},\n    @{old = \"741\u00d73=\"; new = \"695\u00d77=\"},\n    @{old = \"262\u00d77=\"; new = \"619\u00d73=\"}\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $r.old,      # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $r.new,      # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
